# beispiel.xlsx: add two new "real $var" example columns (G, H) to the
# question table, shrink row 3 (it no longer needs to be so tall now that
# the formula text is shorter), widen column B/F, reset the view back to
# the top-left corner and zoomed out a bit, and move the selection to E3.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- new header cells (row 1): G1 = "erg" (same as A1/C1/E1), H1 = "more"
$ws.Range("G1").Value = "erg"
$ws.Range("H1").Value = "more"

# --- new formula/question pair (row 2): G2 = formula, H2 = question text
$ws.Range("G2").Value = "!F*G"
$ws.Range("H2").Value = "Das hier nen test `$F(3,6) mal `$G(2,5)"

# --- column widths: widen B (now longer header wraps) and give the new
# column F (H on-sheet is col 8; col 6 = "F" holds the width bump per the
# cols list) its own width
$ws.Columns.Item(2).ColumnWidth = 43.6
$ws.Columns.Item(6).ColumnWidth = 39.15

# --- row 3 is no longer as tall
$ws.Rows.Item(3).RowHeight = 104.2

# --- view: zoom out, scroll back to the top-left, select E3
$excel.ActiveWindow.Zoom = 140
$excel.ActiveWindow.ScrollRow = 1
$excel.ActiveWindow.ScrollColumn = 1
$ws.Range("E3").Select() | Out-Null
